$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark criteria in C7 and C9 as earned (matching their point value in column B)
$ws.Range("C7").Value = 0.1
$ws.Range("C9").Value = 0.1

# Update the selected cell to reflect where the user ended up after editing
$ws.Range("C10").Select()
